# Generate Report for Handoff
# Replaces the two "handed-back" source records with a freshly generated
# "ready for handoff" report: new source file ids, new handoff timestamps,
# and removal of the old "Latest Target File" / "Latest Handback File"
# columns (F/G) on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$newId1 = "aac1a086-635d-4bc7-8d61-23c24bef0bb6"
$newId2 = "ffff0dc19648-959a-4389-8364-2c5b86ef7bed"

$newMd1 = "$newId1.md"
$newMd2 = "$newId2.md"

$newStatus = "Ready for handoff"
$newOverviewDate = "2016-03-24 05:10:50"
$newHandback = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# NOTE: Range(...).Hyperlinks.Delete() clears the *whole sheet's* hyperlink
# collection in this engine (not just the target range), so remove them all
# up front and rebuild every link once the cell values are in place.
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $newMd1
$ws.Range("B2").Value = $newStatus
$ws.Range("C2").Value = $newStatus
$ws.Range("D2").Value = $newOverviewDate

$ws.Range("A3").Value = $newMd2
$ws.Range("B3").Value = $newStatus
$ws.Range("C3").Value = $newStatus
$ws.Range("D3").Value = $newOverviewDate

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a5d79e98309234ba6b8d61e6236ea74f1a41bf01/e2e/$newMd1", $null, $null, $newMd1) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a5d79e98309234ba6b8d61e6236ea74f1a41bf01/e2e/$newMd2", $null, $null, $newMd2) | Out-Null

# ---------------------------------------------------------------------
# Per-locale sheets (zh-cn, de-de)
# ---------------------------------------------------------------------
$locales = @{
    "zh-cn" = @{ Xlf = "$newId1.c31798f111cbfeadbf8bd1000277a568068217a3.zh-cn.xlf"; HandoffE = "2016-03-24 05:10:46" };
    "de-de" = @{ Xlf = "$newId1.c31798f111cbfeadbf8bd1000277a568068217a3.de-de.xlf"; HandoffE = "2016-03-24 05:10:50" }
}

foreach ($locale in $locales.Keys) {
    $info = $locales[$locale]
    $ws = $wb.Worksheets.Item($locale)

    # Clear every existing hyperlink on the sheet (A2,D2,F2,G2,A3,D3,F3,G3)
    # before touching cell contents.
    $ws.Hyperlinks.Delete()

    # --- Row 2 (was: 362a6b95-... source record) ---
    $ws.Range("A2").Value = $newMd1
    $ws.Range("C2").Value = $newStatus
    $ws.Range("D2").Value = $info.Xlf
    $ws.Range("E2").Value = $info.HandoffE
    $ws.Range("F2").Clear()
    $ws.Range("G2").Clear()
    $ws.Range("H2").Value = $newHandback

    # --- Row 3 (was: 4d150459-... source record) ---
    $ws.Range("A3").Value = $newMd2
    $ws.Range("C3").Value = $newStatus
    $ws.Range("D3").Value = $info.Xlf
    $ws.Range("E3").Value = $info.HandoffE
    $ws.Range("F3").Clear()
    $ws.Range("G3").Clear()
    $ws.Range("H3").Value = $newHandback

    # Rebuild the four hyperlinks that remain (A2, D2, A3, D3) in order so
    # the generated relationship ids line up (rId2..rId5).
    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a5d79e98309234ba6b8d61e6236ea74f1a41bf01/e2e/$newMd1", $null, $null, $newMd1) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3e3d49ae47060ee70802a5ffba08e66f3c2e8db/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/$($info.Xlf)", $null, $null, $info.Xlf) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a5d79e98309234ba6b8d61e6236ea74f1a41bf01/e2e/$newMd2", $null, $null, $newMd2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3e3d49ae47060ee70802a5ffba08e66f3c2e8db/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/ci/ht/$($info.Xlf)", $null, $null, $info.Xlf) | Out-Null
}
